$d = $word.ActiveDocument

function Get-ParaIndexByText([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) { return $i }
    }
    return $null
}

# Replaces the visible text of a plain (no rPr) list-bullet paragraph while preserving
# the document's existing leading empty "<w:r/>" run (Find/Range.Text would merge it away).
function Set-PlainParagraphText([int]$index, [string]$newText) {
    $r = $d.Paragraphs($index).Range
    $inner = $d.Range($r.Start, $r.End - 1)
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $needsSpace = ($newText -ne $newText.Trim())
    $tOpen = "<w:t>"
    if ($needsSpace) { $tOpen = '<w:t xml:space="preserve">' }
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $tOpen + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $inner.InsertXML($xml) | Out-Null
}

# 1. Title change (appears twice: Heading1 at top, and bold recap line near the end).
#    wdFindContinue (1) + wdReplaceAll (2) walks the whole story and replaces every match.
$d.Content.Find.Execute(
    "Play Magic Wilds Free and Uncover the Mystery", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Magic Wilds Free - Exciting Gameplay Features", 2) | Out-Null

# 2. "What we like" bullet list rework.
#    Item 1 text changes, and a brand-new bullet is inserted right after it.
$idx = Get-ParaIndexByText "High winning potential with the Astral Projection feature"
Set-PlainParagraphText $idx "Mysterious magician as the Wild symbol"

# Add a fresh List Bullet paragraph right after it (inherits the ListBullet/spacing/indent pPr).
$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$newBulletXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Astral Projection feature increases winning potential</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs($idx + 1).Range.InsertXML($newBulletXml) | Out-Null

# Remove the old "Magician as the Wild symbol" bullet paragraph entirely.
$removeIdx = Get-ParaIndexByText "The Magician as the Wild symbol increases overall payouts"
$d.Paragraphs($removeIdx).Range.Delete() | Out-Null

# Last "what we like" bullet gets reworded.
$idx = Get-ParaIndexByText "Betting range accommodates all types of players"
Set-PlainParagraphText $idx "Wide betting range and average RTP percentage"

# 3. "What we don't like" bullets reworded.
$idx = Get-ParaIndexByText "RTP percentage is average compared to other online slots"
Set-PlainParagraphText $idx "Limited number of regular symbols"

$idx = Get-ParaIndexByText "Limited variety of symbols on the reels"
Set-PlainParagraphText $idx "No progressive jackpot feature"

# 4. Closing meta blurb (bold title already handled above; italic summary line below).
$d.Content.Find.Execute(
    "Discover the magic of Magic Wilds and play for free. Get the full review on gameplay features, pros & cons, and betting range. Start playing now!",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review to learn more about the exciting gameplay features of Magic Wilds and play free!", 2) | Out-Null
